# Generate Report for Handback
#
# The workbook currently lists two handed-back files (475b1cd4... and
# 551fdd86...) on each of the three sheets (Overview, zh-cn, de-de).
# This run of the report only produced a handback for 475b1cd4..., so the
# 551fdd86... row (row 3) is removed from every sheet, and the surviving
# row's handoff/handback timestamps are refreshed.

$wb = $excel.ActiveWorkbook

function Remove-HyperlinksInRow($ws, $row) {
    # Hyperlinks collection here is a plain, index-backed list: deleting an
    # item while a `foreach` is mid-walk reshuffles later indices and skips
    # entries. Rescan from the top and bail out after each single delete so
    # every match in the target row is actually removed.
    $more = $true
    while ($more) {
        $more = $false
        foreach ($h in $ws.Hyperlinks) {
            if ($h.Range.Row -eq $row) {
                $h.Delete()
                $more = $true
                break
            }
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet: drop the 551fdd86... row (A3:C3) and its hyperlink.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
Remove-HyperlinksInRow $wsOverview 3
$wsOverview.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# zh-cn sheet: drop row 3 + its hyperlinks, refresh row 2's timestamps.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Remove-HyperlinksInRow $wsZhCn 3
$wsZhCn.Rows.Item(3).Delete()
$wsZhCn.Range("E2").Value = "2016-03-11 22:33:05"
$wsZhCn.Range("H2").Value = "2016-03-11 22:33:52"

# ---------------------------------------------------------------------
# de-de sheet: drop row 3 + its hyperlinks, refresh row 2's timestamps.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
Remove-HyperlinksInRow $wsDeDe 3
$wsDeDe.Rows.Item(3).Delete()
$wsDeDe.Range("E2").Value = "2016-03-11 22:33:09"
$wsDeDe.Range("H2").Value = "2016-03-11 22:33:57"
